# Apply cryptos list update (price/volume refresh) per commit on
# Sat Aug 12 02:39:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.392.02'
$ws.Range('D3').Value = '1.848.14'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.15'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07641'
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2934'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.57'
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.00001126'
$ws.Range('E12').Value = '  +12.86%  '
$ws.Range('D13').Value = '1.860.22'
$ws.Range('E13').Value = '  -6.29%  '
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6787'
$ws.Range('E15').Value = '  -0.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.75'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '2.111.23'
$ws.Range('E17').Value = '  -6.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.180'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('D19').Value = '29.409.46'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '228.90'
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.46'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.499'
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.46'
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1398'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.346'
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.465'
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('E30').Value = '  +3.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05599'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.116'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.033'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.853'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7104'
$ws.Range('E36').Value = '  -0.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.586'
$ws.Range('E37').Value = '  -0.48%  '
$ws.Range('D38').Value = '1.240.85'
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01805'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.778'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.414'
$ws.Range('E41').Value = '  +5.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9029'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9999'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.00'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.156'
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.686'
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.989'
$ws.Range('E50').Value = '  -1.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1121'
$ws.Range('E51').Value = '  -0.54%  '
